# Update "want-to-go" counts (column F) on the "展览" (Exhibition),
# "演出" (Performance) and "全部类型" (All types) sheets to the refreshed
# values captured at the newer site-generation run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 255
$ws.Range("F3").Value = 257
$ws.Range("F4").Value = 271
$ws.Range("F5").Value = 2870
$ws.Range("F8").Value = 2227
$ws.Range("F9").Value = 375
$ws.Range("F11").Value = 434
$ws.Range("F13").Value = 2553
$ws.Range("F15").Value = 1356
$ws.Range("F16").Value = 4721
$ws.Range("F18").Value = 5160
$ws.Range("F19").Value = 1722
$ws.Range("F20").Value = 2886
$ws.Range("F21").Value = 3294
$ws.Range("F22").Value = 171
$ws.Range("F23").Value = 1569
$ws.Range("F27").Value = 297
$ws.Range("F28").Value = 1022
$ws.Range("F29").Value = 1910
$ws.Range("F32").Value = 725
$ws.Range("F35").Value = 421

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 22
$ws.Range("F14").Value = 9
$ws.Range("F15").Value = 44

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 255
$ws.Range("F8").Value = 257
$ws.Range("F10").Value = 271
$ws.Range("F11").Value = 2870
$ws.Range("F13").Value = 2227
$ws.Range("F14").Value = 375
$ws.Range("F17").Value = 434
$ws.Range("F19").Value = 22
$ws.Range("F20").Value = 2553
$ws.Range("F21").Value = 1356
$ws.Range("F25").Value = 4721
$ws.Range("F27").Value = 5160
$ws.Range("F28").Value = 1722
$ws.Range("F29").Value = 2886
$ws.Range("F30").Value = 3294
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 171
$ws.Range("F33").Value = 44
$ws.Range("F35").Value = 1569
$ws.Range("F40").Value = 297
$ws.Range("F41").Value = 1022
$ws.Range("F43").Value = 1910
$ws.Range("F46").Value = 725
$ws.Range("F49").Value = 421
